$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2024-10-02 15:16:26"

$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "2024-10-05 09:50:46"

$ws.Range("A4").Value = "test hello"
$ws.Range("B4").Value = "2024-10-05 10:30:25"
